$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The author added two new bullet points ("製作kernal" and
# "將kernal與image運算後得出新的list") right before the existing
# "根據狀況判斷哪個pixel需要上色" bullet. Everything else in the list keeps
# its original text/formatting untouched, so the cleanest way to reproduce
# the edit is to insert two brand-new paragraphs in front of that bullet,
# using InsertXML so the exact run layout (incl. <w:proofErr/> spell-check
# markers around "kernal") matches what Word itself would produce.
# ---------------------------------------------------------------------------

function New-WordXmlFragment([string]$InnerParagraphXml) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        $InnerParagraphXml +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Shared paragraph properties used by every bullet in this list (pStyle a3 /
# numbered list numId 1 / ind leftChars 0 / 標楷體 36pt paragraph mark).
$pPr = '<w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr>'
$rPr = '<w:rPr><w:rFonts w:ascii="標楷體" w:eastAsia="標楷體" w:hAnsi="標楷體" w:hint="eastAsia"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr>'

# "製作" + spell-checked "kernal"
$p1Inner = "<w:p>$pPr<w:r>$rPr<w:t>製作</w:t></w:r>" +
    '<w:proofErr w:type="spellStart"/>' +
    "<w:r>$rPr<w:t>kernal</w:t></w:r>" +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'

# "將" + spell-checked "kernal" + "與image運算後得出新的list"
$p2Inner = "<w:p>$pPr<w:r>$rPr<w:t>將</w:t></w:r>" +
    '<w:proofErr w:type="spellStart"/>' +
    "<w:r>$rPr<w:t>kernal</w:t></w:r>" +
    '<w:proofErr w:type="spellEnd"/>' +
    "<w:r>$rPr<w:t>與image運算後得出新的list</w:t></w:r>" +
    '</w:p>'

$frag1 = New-WordXmlFragment $p1Inner
$frag2 = New-WordXmlFragment $p2Inner

# Find the existing bullet the two new ones get inserted in front of, and
# resolve it to a plain $d.Paragraphs index (operating through
# $d.Paragraphs.Item(n) keeps paragraph navigation such as .Previous()
# reliable, unlike paragraphs reached via a Find range).
$find = $d.Content
$find.Find.Execute("根據狀況判斷哪個pixel需要上色", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$anchorIndex = $find.Paragraphs.Item(1).Index

# Insert an empty paragraph above the anchor, then fill it with the first
# new bullet's real OOXML (runs + proofErr markers).
$d.Paragraphs.Item($anchorIndex).Range.InsertParagraphBefore() | Out-Null
$d.Paragraphs.Item($anchorIndex).Range.InsertXML($frag1) | Out-Null

# The anchor bullet is now one slot further down; repeat for the second
# new bullet, still directly above the original "根據狀況..." paragraph.
$anchorIndex = $anchorIndex + 1
$d.Paragraphs.Item($anchorIndex).Range.InsertParagraphBefore() | Out-Null
$d.Paragraphs.Item($anchorIndex).Range.InsertXML($frag2) | Out-Null

Write-Output "done"
